$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.791.43'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.687.66'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '647.79'
$ws.Range("E5").Value = '  -4.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.26'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.503'
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.146'
$ws.Range("E9").Value = '  -0.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.18'
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.444'
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.311.20'
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.77'
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.692.01'
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.848.18'
$ws.Range("E16").Value = '  +0.49%  '
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.05'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.40'
$ws.Range("E20").Value = '  +5.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '472.09'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.652'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '80.13'
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.837.72'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.98'
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.15'
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.65'
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.71'
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.01'
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.54'
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.84'
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.165'
$ws.Range("E35").Value = '  +1.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.686.90'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.43'
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.89'
$ws.Range("E39").Value = '  -5.29%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '178.67'
$ws.Range("E40").Value = '  +6.98%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.25'
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0904'
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.931'
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.83'
$ws.Range("E45").Value = '  +2.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.84'
$ws.Range("E46").Value = '  +0.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '29.23'
$ws.Range("E47").Value = '  +4.12%  '
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.26'
$ws.Range("E49").Value = '  -2.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.85'
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("E51").Value = '  -3.70%  '
